$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.681069641154807
$ws.Range("C2").Value = 0.1579971622521441
$ws.Range("D2").Value = 0.08123152478877316
$ws.Range("F2").Value = 2.522667468636485
$ws.Range("G2").Value = 1.848888428779702
$ws.Range("H2").Value = 1.591304099586381
$ws.Range("J2").Value = 0.2043006506243119
$ws.Range("L2").Value = 0.3189709388046893
$ws.Range("M2").Value = 0.3855908463116364

$ws.Range("B3").Value = 1.594119468384292
$ws.Range("C3").Value = 0.1434501723039716
$ws.Range("D3").Value = 0.08110058299276091
$ws.Range("F3").Value = 2.524891366810763
$ws.Range("G3").Value = 1.845586416402725
$ws.Range("H3").Value = 1.596441010221639
$ws.Range("J3").Value = 0.2054991594563775
$ws.Range("L3").Value = 0.3167682993302279
$ws.Range("M3").Value = 0.3725014846121795

$ws.Range("B4").Value = 1.541434001318976
$ws.Range("C4").Value = 0.1344554492388568
$ws.Range("D4").Value = 0.08102483997508791
$ws.Range("F4").Value = 2.527626293406257
$ws.Range("G4").Value = 1.844728067973762
$ws.Range("H4").Value = 1.600380888814399
$ws.Range("J4").Value = 0.206272965773836
$ws.Range("L4").Value = 0.3155293639560384
$ws.Range("M4").Value = 0.3646297798263305

$ws.Range("B5").Value = 1.520141633494234
$ws.Range("C5").Value = 0.1307742864037493
$ws.Range("D5").Value = 0.08099514889536152
$ws.Range("F5").Value = 2.529084715873367
$ws.Range("G5").Value = 1.844671567360677
$ws.Range("H5").Value = 1.602183841260469
$ws.Range("J5").Value = 0.2065978500984869
$ws.Range("L5").Value = 0.3150531084479837
$ws.Range("M5").Value = 0.3614636982544894

$ws.Range("B6").Value = 1.516616789578961
$ws.Range("C6").Value = 0.1301620831302444
$ws.Range("D6").Value = 0.08099028979523659
$ws.Range("F6").Value = 2.529347644767171
$ws.Range("G6").Value = 1.844679878718054
$ws.Range("H6").Value = 1.602495138744189
$ws.Range("J6").Value = 0.2066523743214672
$ws.Range("L6").Value = 0.3149757576295897
$ws.Range("M6").Value = 0.3609404962213389

$ws.Range("B7").Value = 1.541146125702596
$ws.Range("C7").Value = 0.1344058674453379
$ws.Range("D7").Value = 0.0810244347872171
$ws.Range("F7").Value = 2.527644570286782
$ws.Range("G7").Value = 1.844726119424635
$ws.Range("H7").Value = 1.600404404956706
$ws.Range("J7").Value = 0.206277308579061
$ws.Range("L7").Value = 0.3155228250140283
$ws.Range("M7").Value = 0.3645869118826681

$ws.Range("B8").Value = 1.65094389943738
$ws.Range("C8").Value = 0.1529944245857848
$ws.Range("D8").Value = 0.08118541206290786
$ws.Range("F8").Value = 2.523149769451521
$ws.Range("G8").Value = 1.847506822801705
$ws.Range("H8").Value = 1.592912106978872
$ws.Range("J8").Value = 0.2047060373486898
$ws.Range("L8").Value = 0.3181879603897713
$ws.Range("M8").Value = 0.3810434321109071

$ws.Range("B9").Value = 1.871808056873363
$ws.Range("C9").Value = 0.1889475467847035
$ws.Range("D9").Value = 0.08153788961709907
$ws.Range("F9").Value = 2.525224125302842
$ws.Range("G9").Value = 1.86227029203539
$ws.Range("H9").Value = 1.584463777751751
$ws.Range("J9").Value = 0.2019248151364299
$ws.Range("L9").Value = 0.3243118086517995
$ws.Range("M9").Value = 0.4146208124806279

$ws.Range("B10").Value = 2.037452291848524
$ws.Range("C10").Value = 0.2150599816505974
$ws.Range("D10").Value = 0.08181916066350503
$ws.Range("F10").Value = 2.533421776865836
$ws.Range("G10").Value = 1.878844412572391
$ws.Range("H10").Value = 1.582077754573959
$ws.Range("J10").Value = 0.2000631806808641
$ws.Range("L10").Value = 0.3293549951859376
$ws.Range("M10").Value = 0.4400832799449717

$ws.Range("B11").Value = 2.113540803863771
$ws.Range("C11").Value = 0.2268741657499334
$ws.Range("D11").Value = 0.08195194215216617
$ws.Range("F11").Value = 2.5386082011361
$ws.Range("G11").Value = 1.887639422631281
$ws.Range("H11").Value = 1.581825242158772
$ws.Range("J11").Value = 0.1992554840137868
$ws.Range("L11").Value = 0.3317667993658375
$ws.Range("M11").Value = 0.4518386420009648

$ws.Range("B12").Value = 2.14245898436684
$ws.Range("C12").Value = 0.231338611864544
$ws.Range("D12").Value = 0.08200291545768401
$ws.Range("F12").Value = 2.540782298960437
$ws.Range("G12").Value = 1.891151230368081
$ws.Range("H12").Value = 1.581849633759845
$ws.Range("J12").Value = 0.198955245011768
$ws.Range("L12").Value = 0.3326969379065048
$ws.Range("M12").Value = 0.4563147811982518

$ws.Range("B13").Value = 2.13622627344796
$ws.Range("C13").Value = 0.2303775292680257
$ws.Range("D13").Value = 0.08199190671459711
$ws.Range("F13").Value = 2.540304714413409
$ws.Range("G13").Value = 1.890386822201663
$ws.Range("H13").Value = 1.581839039117227
$ws.Range("J13").Value = 0.1990196572101111
$ws.Range("L13").Value = 0.33249586811732
$ws.Range("M13").Value = 0.4553496710089462

$ws.Range("B14").Value = 2.115917817536854
$ws.Range("C14").Value = 0.2272416456886788
$ws.Range("D14").Value = 0.08195612190568546
$ws.Range("F14").Value = 2.53878285132086
$ws.Range("G14").Value = 1.887924702435697
$ws.Range("H14").Value = 1.581824842245027
$ws.Range("J14").Value = 0.1992306706555658
$ws.Range("L14").Value = 0.3318429853363796
$ws.Range("M14").Value = 0.4522064036944755

$ws.Range("B15").Value = 2.103491969008928
$ws.Range("C15").Value = 0.225319609960934
$ws.Range("D15").Value = 0.08193429269901564
$ws.Range("F15").Value = 2.537878045554692
$ws.Range("G15").Value = 1.886440221336358
$ws.Range("H15").Value = 1.581831782656394
$ws.Range("J15").Value = 0.1993606537133639
$ws.Range("L15").Value = 0.3314452667282808
$ws.Range("M15").Value = 0.4502842661160571

$ws.Range("B16").Value = 2.032494477435705
$ws.Range("C16").Value = 0.2142865943992831
$ws.Range("D16").Value = 0.08181058003443731
$ws.Range("F16").Value = 2.533112203471433
$ws.Range("G16").Value = 1.878294968139585
$ws.Range("H16").Value = 1.582111037631336
$ws.Range("J16").Value = 0.2001167525240044
$ws.Range("L16").Value = 0.3291997387648138
$ws.Range("M16").Value = 0.4393184950949873

$ws.Range("B17").Value = 1.989127819886903
$ws.Range("C17").Value = 0.2075016383611228
$ws.Range("D17").Value = 0.08173592136945729
$ws.Range("F17").Value = 2.530562151726315
$ws.Range("G17").Value = 1.873620221810882
$ws.Range("H17").Value = 1.582495849343275
$ws.Range("J17").Value = 0.2005906167548499
$ws.Range("L17").Value = 0.3278522551425453
$ws.Range("M17").Value = 0.4326353896190795

$ws.Range("B18").Value = 1.964253795383513
$ws.Range("C18").Value = 0.2035930437600371
$ws.Range("D18").Value = 0.0816934344361826
$ws.Range("F18").Value = 2.52923256122321
$ws.Range("G18").Value = 1.871049543276001
$ws.Range("H18").Value = 1.582795564905098
$ws.Range("J18").Value = 0.2008668586624811
$ws.Range("L18").Value = 0.3270882948112188
$ws.Range("M18").Value = 0.4288076757635295

$ws.Range("B19").Value = 1.955843808474071
$ws.Range("C19").Value = 0.2022686204390709
$ws.Range("D19").Value = 0.08167912728123738
$ws.Range("F19").Value = 2.528805920285734
$ws.Range("G19").Value = 1.870199416897606
$ws.Range("H19").Value = 1.582910498274288
$ws.Range("J19").Value = 0.2009610232401249
$ws.Range("L19").Value = 0.3268315356893225
$ws.Range("M19").Value = 0.4275144699343016

$ws.Range("B20").Value = 1.993737104670856
$ws.Range("C20").Value = 0.208224536863213
$ws.Range("D20").Value = 0.08174382186598805
$ws.Range("F20").Value = 2.530819412471999
$ws.Range("G20").Value = 1.874105626577233
$ws.Range("H20").Value = 1.58244677151302
$ws.Range("J20").Value = 0.200539791537957
$ws.Range("L20").Value = 0.3279945512339424
$ws.Range("M20").Value = 0.4333451387629808

$ws.Range("B21").Value = 2.121880056400016
$ws.Range("C21").Value = 0.2281629843205053
$ws.Range("D21").Value = 0.08196661401516891
$ws.Range("F21").Value = 2.539224152575784
$ws.Range("G21").Value = 1.888642959035082
$ws.Range("H21").Value = 1.581825753161326
$ws.Range("J21").Value = 0.1991685385293644
$ws.Range("L21").Value = 0.3320342963762499
$ws.Range("M21").Value = 0.4531289899579676

$ws.Range("B22").Value = 2.206241229643183
$ws.Range("C22").Value = 0.2411395124177318
$ws.Range("D22").Value = 0.08211625306283388
$ws.Range("F22").Value = 2.545942087779608
$ws.Range("G22").Value = 1.8992012434922
$ws.Range("H22").Value = 1.582119520861795
$ws.Range("J22").Value = 0.1983050859971454
$ws.Range("L22").Value = 0.3347726296016447
$ws.Range("M22").Value = 0.4662024283841291

$ws.Range("B23").Value = 2.161160293229216
$ws.Range("C23").Value = 0.2342186901144316
$ws.Range("D23").Value = 0.0820360198250345
$ws.Range("F23").Value = 2.542244331139941
$ws.Range("G23").Value = 1.893469077649939
$ws.Range("H23").Value = 1.581898633884265
$ws.Range("J23").Value = 0.198762935820314
$ws.Range("L23").Value = 0.3333021754959304
$ws.Range("M23").Value = 0.4592118009536819

$ws.Range("B24").Value = 1.991653066164588
$ws.Range("C24").Value = 0.2078977388457872
$ws.Range("D24").Value = 0.0817402486961285
$ws.Range("F24").Value = 2.530702679845362
$ws.Range("G24").Value = 1.873885811026838
$ws.Range("H24").Value = 1.582468715150725
$ws.Range("J24").Value = 0.2005627577462317
$ws.Range("L24").Value = 0.3279301857562587
$ws.Range("M24").Value = 0.4330242160222326

$ws.Range("B25").Value = 1.811465155173494
$ws.Range("C25").Value = 0.1792745122400277
$ws.Range("D25").Value = 0.08143860959734184
$ws.Range("F25").Value = 2.523493581213145
$ws.Range("G25").Value = 1.857274315135044
$ws.Range("H25").Value = 1.586079322376548
$ws.Range("J25").Value = 0.2026452112190757
$ws.Range("L25").Value = 0.3225593903743444
$ws.Range("M25").Value = 0.4053977410139069
